$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Credit Limit (C) and Current Outstanding (D) figures
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 5

$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 6

$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 7

$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 8

# Re-format the "Date Data" column with a custom date format
$ws.Range("E2:E5").NumberFormat = "dd\-mmm\-yyyy"

# Move the active selection
$ws.Range("G6").Select() | Out-Null
